$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wv = $wb.Windows.Item(1)
Get-Member -InputObject $wv
